$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 72.266001
$ws.Range("H2").Value = 216.798003
$ws.Range("I2").Value = 0.2949652269937106
$ws.Range("J2").Value = 0.2949652269937106
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09554499999999999
$ws.Range("N2").Value = 0.286635
$ws.Range("O2").Value = 0.01029975823317688
$ws.Range("P2").Value = 0.01029975823317688
$ws.Range("Q2").Value = 6.904655065545
$ws.Range("R2").Value = 62.14189558990499
$ws.Range("S2").Value = 0.003038070525229358
$ws.Range("T2").Value = 0.003038070525229358

# Row 3
$ws.Range("G3").Value = 72.266001
$ws.Range("H3").Value = 216.798003
$ws.Range("I3").Value = 0.2949652269937106
$ws.Range("J3").Value = 0.2949652269937106
$ws.Range("O3").Value = 0.4011437372432085
$ws.Range("P3").Value = 0.4011437372432086
$ws.Range("Q3").Value = 268.914966221849
$ws.Range("R3").Value = 2420.234695996641
$ws.Range("S3").Value = 0.1183234535130484
$ws.Range("T3").Value = 0.1183234535130484

# Row 4
$ws.Range("G4").Value = 72.266001
$ws.Range("H4").Value = 216.798003
$ws.Range("I4").Value = 0.2949652269937106
$ws.Range("J4").Value = 0.2949652269937106
$ws.Range("M4").Value = 5.459703999999999
$ws.Range("N4").Value = 16.379112
$ws.Range("O4").Value = 0.5885565045236145
$ws.Range("P4").Value = 0.5885565045236146
$ws.Range("Q4").Value = 394.550974723704
$ws.Range("R4").Value = 3550.958772513336
$ws.Range("S4").Value = 0.1736037029554328
$ws.Range("T4").Value = 0.1736037029554328

# Row 5
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5578585839920717
$ws.Range("J5").Value = 0.5578585839920718
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09554499999999999
$ws.Range("N5").Value = 0.286635
$ws.Range("O5").Value = 0.01029975823317688
$ws.Range("P5").Value = 0.01029975823317688
$ws.Range("Q5").Value = 13.05855994307
$ws.Range("R5").Value = 117.52703948763
$ws.Range("S5").Value = 0.005745808543420737
$ws.Range("T5").Value = 0.005745808543420738

# Row 6
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5578585839920717
$ws.Range("J6").Value = 0.5578585839920718
$ws.Range("O6").Value = 0.4011437372432085
$ws.Range("P6").Value = 0.4011437372432086
$ws.Range("Q6").Value = 508.5905338733206
$ws.Range("R6").Value = 4577.314804859885
$ws.Range("S6").Value = 0.223781477235784
$ws.Range("T6").Value = 0.223781477235784

# Row 7
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5578585839920717
$ws.Range("J7").Value = 0.5578585839920718
$ws.Range("M7").Value = 5.459703999999999
$ws.Range("N7").Value = 16.379112
$ws.Range("O7").Value = 0.5885565045236145
$ws.Range("P7").Value = 0.5885565045236146
$ws.Range("Q7").Value = 746.2020195239838
$ws.Range("R7").Value = 6715.818175715855
$ws.Range("S7").Value = 0.3283312982128669
$ws.Range("T7").Value = 0.3283312982128671

# Row 8
$ws.Range("G8").Value = 36.057927
$ws.Range("H8").Value = 108.173781
$ws.Range("I8").Value = 0.1471761890142177
$ws.Range("J8").Value = 0.1471761890142177
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.09554499999999999
$ws.Range("N8").Value = 0.286635
$ws.Range("O8").Value = 0.01029975823317688
$ws.Range("P8").Value = 0.01029975823317688
$ws.Range("Q8").Value = 3.445154635215
$ws.Range("R8").Value = 31.006391716935
$ws.Range("S8").Value = 0.001515879164526786
$ws.Range("T8").Value = 0.001515879164526786

# Row 9
$ws.Range("G9").Value = 36.057927
$ws.Range("H9").Value = 108.173781
$ws.Range("I9").Value = 0.1471761890142177
$ws.Range("J9").Value = 0.1471761890142177
$ws.Range("O9").Value = 0.4011437372432085
$ws.Range("P9").Value = 0.4011437372432086
$ws.Range("Q9").Value = 134.178120929023
$ws.Range("R9").Value = 1207.603088361207
$ws.Range("S9").Value = 0.05903880649437614
$ws.Range("T9").Value = 0.05903880649437614

# Row 10
$ws.Range("G10").Value = 36.057927
$ws.Range("H10").Value = 108.173781
$ws.Range("I10").Value = 0.1471761890142177
$ws.Range("J10").Value = 0.1471761890142177
$ws.Range("M10").Value = 5.459703999999999
$ws.Range("N10").Value = 16.379112
$ws.Range("O10").Value = 0.5885565045236145
$ws.Range("P10").Value = 0.5885565045236146
$ws.Range("Q10").Value = 196.865608273608
$ws.Range("R10").Value = 1771.790474462472
$ws.Range("S10").Value = 0.08662150335531477
$ws.Range("T10").Value = 0.08662150335531478
